# Apply the mapping update described by the commit:
#  - Update the Metadata "Date" value to the new timestamp
#  - Correct the Source value in "Mapping Table 1" row 5 (was pointing to a
#    duplicated/incorrect "FRCDADICOMExamenImagerie.description" string,
#    should point to "FRCDADICOMExamenImagerie.text" like the corresponding
#    target in "Mapping Table 0")

$wb = $excel.ActiveWorkbook

# --- Update Date on the Metadata sheet ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# --- Fix the Source mapping entry on "Mapping Table 1" ---
$mapWs = $wb.Worksheets.Item("Mapping Table 1")
$mapWs.Range("A5").Value = "FRCDADICOMExamenImagerie.text"
